$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction16")

# Clear the row entirely first so stale cells (C1:S1) are removed
$ws.Rows.Item(1).Clear()

# Set the new random, non-overlapping values
$ws.Range("A1").Value = 32
$ws.Range("B1").Value = 33
